# Commit: "ok - Commented GitLatch Commit @ 2024-6-25-7-13-45-70"
#
# The document has a single paragraph containing one run of text
# "Ewc123122333www". The edit prepends a new run containing the text
# "ok" (same run formatting: lang="en-IN") directly in front of the
# existing run, inside the very same paragraph - i.e. two sibling
# <w:r> elements end up inside the original <w:p>.
#
# A plain Range.InsertBefore("ok") would work text-wise, but since the
# inserted text shares identical run formatting with the text that
# follows it, the engine coalesces it into a single <w:r> run instead
# of leaving two separate runs behind. To keep "ok" as its own
# <w:r>/<w:t> run (matching the target OOXML exactly) we insert it as
# a tiny WordprocessingML package via Range.InsertXML - that always
# lands in its own run - which creates it as a new leading paragraph,
# then we delete the paragraph mark that separates it from the
# original paragraph so the two runs merge back into one <w:p>.

$d = $word.ActiveDocument

# Collapsed range at the very start of the document/paragraph.
$insertionPoint = $d.Range(0, 0)

$openXmlPackage = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
              <w:t>ok</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$null = $insertionPoint.InsertXML($openXmlPackage)

# InsertXML dropped the "ok" run in as its own leading paragraph, so
# the document currently reads: [P1: "ok"][P2: "Ewc123122333www"].
# Remove the paragraph mark that ends P1 to fold its run into P2,
# leaving a single paragraph with the two runs as siblings, exactly
# like the target diff.
$firstParagraph = $d.Paragraphs(1)
$paragraphMark = $d.Range($firstParagraph.Range.End - 1, $firstParagraph.Range.End)
$paragraphMark.Delete()
